$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articulo")
$ws.Activate()
$ws.Range("A2").Value = "xbox"
$ws.Range("B7").Select()
